# Restored from revision of admin on 01/11/2021 12:26:52 PM.TEST Author: admin. Type: SAVE.
# Change cell C10 on the active worksheet from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
